# Reorders the comma-separated "Recorded By" values in column G.
# Rule (derived from target diff):
#   - If the token "System" is present in the list, move it to the front
#     of the list while preserving the relative order of the remaining
#     tokens.
#   - Otherwise (no "System" token present), reverse the order of the
#     tokens.
#   - Lists with a single token are left unchanged (no-op either way).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "

    if ($parts.Count -le 1) { continue }

    $idx = [array]::IndexOf($parts, "System")

    if ($idx -ge 0) {
        $rest = @()
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $idx) {
                $rest += $parts[$i]
            }
        }
        $newParts = @("System") + $rest
    }
    else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
